$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 647, pushing the existing row 647 (and everything
# below it) down by one row.
$ws.Rows.Item(647).Insert()

# Populate the newly-inserted row 647 with the new weekly record.
$ws.Cells.Item(647, 1).Value = 10
$ws.Cells.Item(647, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(647, 3).Value = "La Araucanía"
$ws.Cells.Item(647, 4).Value = 45265
$ws.Cells.Item(647, 5).Value = 9
$ws.Cells.Item(647, 6).Value = 100112040
$ws.Cells.Item(647, 7).Value = "Cilantro"
$ws.Cells.Item(647, 8).Value = "Sin especificar"
$ws.Cells.Item(647, 9).Value = "Primera"
$ws.Cells.Item(647, 10).Value = 30
$ws.Cells.Item(647, 11).Value = 7000
$ws.Cells.Item(647, 12).Value = 7000
$ws.Cells.Item(647, 13).Value = 7000
$ws.Cells.Item(647, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(647, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(647, 16).Value = 3500
$ws.Cells.Item(647, 17).Value = 2
$ws.Cells.Item(647, 18).Value = "Hortaliza"
